# Remove leftover template values (OCD-4438)
# The "Template" sheet shipped with two stray formula results (U5, U6)
# left over from a previous report run. Clear them out but keep the
# existing cell formatting (style s="3") intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
